# Scheduled runner update: refresh cached market-board price/profit
# columns (H-N) across several recipe sheets.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 12
$wsALC.Range("H12").Value = 100
$wsALC.Range("I12").Value = 100
$wsALC.Range("J12").Value = 0
$wsALC.Range("K12").Value = 100
$wsALC.Range("L12").Value = 0
$wsALC.Range("M12").Value = 70
$wsALC.Range("N12").ClearContents()

# ALC row 33
$wsALC.Range("H33").Value = 656.2857
$wsALC.Range("I33").Value = 835.1429000000001
$wsALC.Range("J33").Value = 298.57144
$wsALC.Range("K33").Value = 835.1429000000001
$wsALC.Range("L33").Value = 298.57144
$wsALC.Range("M33").Value = -606.1429000000001
$wsALC.Range("N33").Value = -756.5714399999999

# ALC row 112
$wsALC.Range("H112").Value = 1063.8889
$wsALC.Range("I112").Value = 547.1429000000001
$wsALC.Range("J112").Value = 2872.5
$wsALC.Range("K112").Value = 1641.4287
$wsALC.Range("L112").Value = 8617.5
$wsALC.Range("M112").Value = -533.4287000000002
$wsALC.Range("N112").Value = -10833.5

# ALC row 136
$wsALC.Range("H136").Value = 85192.71000000001
$wsALC.Range("J136").Value = 85192.71000000001
$wsALC.Range("L136").Value = 85192.71000000001
$wsALC.Range("N136").Value = -95392.71000000001

# ARM row 32
$wsARM.Range("H32").Value = 28418.014
$wsARM.Range("I32").Value = 10790.307
$wsARM.Range("K32").Value = 10790.307
$wsARM.Range("M32").Value = -10503.307

# ARM row 74
$wsARM.Range("H74").Value = 20794.867
$wsARM.Range("I74").Value = 1257.7561
$wsARM.Range("J74").Value = 87546.664
$wsARM.Range("K74").Value = 1257.7561
$wsARM.Range("L74").Value = 87546.664
$wsARM.Range("M74").Value = -383.7561000000001
$wsARM.Range("N74").Value = -89294.664

# ARM row 77
$wsARM.Range("H77").Value = 20794.867
$wsARM.Range("I77").Value = 1257.7561
$wsARM.Range("J77").Value = 87546.664
$wsARM.Range("K77").Value = 6288.780500000001
$wsARM.Range("L77").Value = 437733.32
$wsARM.Range("M77").Value = -1920.780500000001
$wsARM.Range("N77").Value = -446469.32

# ARM row 102
$wsARM.Range("H102").Value = 62501624
$wsARM.Range("I102").Value = 1410.6364
$wsARM.Range("J102").Value = 200002100
$wsARM.Range("K102").Value = 1410.6364
$wsARM.Range("L102").Value = 200002100
$wsARM.Range("M102").Value = 211.3635999999999
$wsARM.Range("N102").Value = -200005344

# ARM row 107
$wsARM.Range("H107").Value = 15876
$wsARM.Range("J107").Value = 15876
$wsARM.Range("L107").Value = 15876
$wsARM.Range("N107").Value = -23556

# BSM row 35
$wsBSM.Range("H35").Value = 33849.715
$wsBSM.Range("J35").Value = 33849.715
$wsBSM.Range("L35").Value = 33849.715
$wsBSM.Range("N35").Value = -34469.715

# BSM row 80
$wsBSM.Range("H80").Value = 257.73914
$wsBSM.Range("I80").Value = 70.85714
$wsBSM.Range("J80").Value = 339.5
$wsBSM.Range("K80").Value = 70.85714
$wsBSM.Range("L80").Value = 339.5
$wsBSM.Range("M80").Value = 927.14286
$wsBSM.Range("N80").Value = -2335.5

# BSM row 82
$wsBSM.Range("H82").Value = 24203.65
$wsBSM.Range("J82").Value = 36889.637
$wsBSM.Range("L82").Value = 36889.637
$wsBSM.Range("N82").Value = -37655.637

# BSM row 83
$wsBSM.Range("H83").Value = 257.73914
$wsBSM.Range("I83").Value = 70.85714
$wsBSM.Range("J83").Value = 339.5
$wsBSM.Range("K83").Value = 354.2857
$wsBSM.Range("L83").Value = 1697.5
$wsBSM.Range("M83").Value = 4637.7143
$wsBSM.Range("N83").Value = -11681.5

# BSM row 85
$wsBSM.Range("H85").Value = 24203.65
$wsBSM.Range("J85").Value = 36889.637
$wsBSM.Range("L85").Value = 36889.637
$wsBSM.Range("N85").Value = -39541.637

# CRP row 7
$wsCRP.Range("H7").Value = 109.82353
$wsCRP.Range("I7").Value = 43
$wsCRP.Range("J7").Value = 156.6
$wsCRP.Range("K7").Value = 43
$wsCRP.Range("L7").Value = 156.6
$wsCRP.Range("M7").Value = 70
$wsCRP.Range("N7").Value = -382.6

# CRP row 16
$wsCRP.Range("H16").Value = 4208167.5
$wsCRP.Range("I16").Value = 7361193.5
$wsCRP.Range("J16").Value = 4133.3335
$wsCRP.Range("K16").Value = 7361193.5
$wsCRP.Range("L16").Value = 4133.3335
$wsCRP.Range("M16").Value = -7360906.5
$wsCRP.Range("N16").Value = -4707.3335

# CRP row 41
$wsCRP.Range("H41").Value = 21072.25
$wsCRP.Range("I41").Value = 9999
$wsCRP.Range("J41").Value = 24763.334
$wsCRP.Range("K41").Value = 9999
$wsCRP.Range("L41").Value = 24763.334
$wsCRP.Range("M41").Value = -9571
$wsCRP.Range("N41").Value = -25619.334

# CRP row 50
$wsCRP.Range("H50").Value = 9192.5
$wsCRP.Range("J50").Value = 9192.5
$wsCRP.Range("L50").Value = 9192.5
$wsCRP.Range("N50").Value = -10442.5

# CRP row 51
$wsCRP.Range("H51").Value = 9234.333000000001
$wsCRP.Range("J51").Value = 9234.333000000001
$wsCRP.Range("L51").Value = 9234.333000000001
$wsCRP.Range("N51").Value = -10706.333

# CRP row 60
$wsCRP.Range("H60").Value = 25603.77
$wsCRP.Range("J60").Value = 25603.77
$wsCRP.Range("L60").Value = 25603.77
$wsCRP.Range("N60").Value = -26625.77

# CRP row 61
$wsCRP.Range("H61").Value = 9234.333000000001
$wsCRP.Range("J61").Value = 9234.333000000001
$wsCRP.Range("L61").Value = 9234.333000000001
$wsCRP.Range("N61").Value = -9930.333000000001

# CRP row 68
$wsCRP.Range("H68").Value = 17399.5
$wsCRP.Range("J68").Value = 17399.5
$wsCRP.Range("L68").Value = 17399.5
$wsCRP.Range("N68").Value = -18897.5

# CRP row 71
$wsCRP.Range("H71").Value = 17399.5
$wsCRP.Range("J71").Value = 17399.5
$wsCRP.Range("L71").Value = 52198.5
$wsCRP.Range("N71").Value = -59686.5

# CRP row 113
$wsCRP.Range("H113").Value = 4208167.5
$wsCRP.Range("I113").Value = 7361193.5
$wsCRP.Range("J113").Value = 4133.3335
$wsCRP.Range("K113").Value = 7361193.5
$wsCRP.Range("L113").Value = 4133.3335
$wsCRP.Range("M113").Value = -7359023.5
$wsCRP.Range("N113").Value = -8473.333500000001

# CRP row 131
$wsCRP.Range("H131").Value = 43435.332
$wsCRP.Range("J131").Value = 43435.332
$wsCRP.Range("L131").Value = 43435.332
$wsCRP.Range("N131").Value = -53515.332

# CUL row 12
$wsCUL.Range("H12").Value = 540.76
$wsCUL.Range("I12").Value = 166.14285
$wsCUL.Range("J12").Value = 686.44446
$wsCUL.Range("K12").Value = 498.42855
$wsCUL.Range("L12").Value = 2059.33338
$wsCUL.Range("M12").Value = -325.42855
$wsCUL.Range("N12").Value = -2405.33338

# CUL row 23
$wsCUL.Range("H23").Value = 63.52381
$wsCUL.Range("I23").Value = 16.25
$wsCUL.Range("J23").Value = 74.64706
$wsCUL.Range("K23").Value = 48.75
$wsCUL.Range("L23").Value = 223.94118
$wsCUL.Range("M23").Value = 186.25
$wsCUL.Range("N23").Value = -693.94118

# CUL row 70
$wsCUL.Range("H70").Value = 6733.7856
$wsCUL.Range("I70").Value = 2943.5
$wsCUL.Range("J70").Value = 8249.9
$wsCUL.Range("K70").Value = 8830.5
$wsCUL.Range("L70").Value = 24749.7
$wsCUL.Range("M70").Value = -8515.5
$wsCUL.Range("N70").Value = -25379.7

# CUL row 73
$wsCUL.Range("H73").Value = 6733.7856
$wsCUL.Range("I73").Value = 2943.5
$wsCUL.Range("J73").Value = 8249.9
$wsCUL.Range("K73").Value = 8830.5
$wsCUL.Range("L73").Value = 24749.7
$wsCUL.Range("M73").Value = -7738.5
$wsCUL.Range("N73").Value = -26933.7

# CUL row 122
$wsCUL.Range("H122").Value = 896.9636
$wsCUL.Range("I122").Value = 410.33334
$wsCUL.Range("J122").Value = 1366.2142
$wsCUL.Range("K122").Value = 3693.00006
$wsCUL.Range("L122").Value = 12295.9278
$wsCUL.Range("M122").Value = -1243.00006
$wsCUL.Range("N122").Value = -17195.9278

# GSM row 12
$wsGSM.Range("H12").Value = 3444.7778
$wsGSM.Range("I12").Value = 1571.8572
$wsGSM.Range("J12").Value = 10000
$wsGSM.Range("K12").Value = 1571.8572
$wsGSM.Range("L12").Value = 10000
$wsGSM.Range("M12").Value = -1431.8572
$wsGSM.Range("N12").Value = -10280

# GSM row 123
$wsGSM.Range("H123").Value = 14965.923
$wsGSM.Range("J123").Value = 14965.923
$wsGSM.Range("L123").Value = 14965.923
$wsGSM.Range("N123").Value = -19865.923

# GSM row 132
$wsGSM.Range("H132").Value = 2523.2222
$wsGSM.Range("I132").Value = 2305.8462
$wsGSM.Range("J132").Value = 3088.4
$wsGSM.Range("K132").Value = 6917.5386
$wsGSM.Range("L132").Value = 9265.200000000001
$wsGSM.Range("M132").Value = -4387.5386
$wsGSM.Range("N132").Value = -14325.2

# WVR row 118
$wsWVR.Range("H118").Value = 29896
$wsWVR.Range("J118").Value = 29896
$wsWVR.Range("L118").Value = 29896
$wsWVR.Range("N118").Value = -33210

# WVR row 136
$wsWVR.Range("H136").Value = 3279.2341
$wsWVR.Range("I136").Value = 3463.279
$wsWVR.Range("J136").Value = 1300.75
$wsWVR.Range("K136").Value = 10389.837
$wsWVR.Range("L136").Value = 3902.25
$wsWVR.Range("M136").Value = -7839.837
$wsWVR.Range("N136").Value = -9002.25
